$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.414.05"
$ws.Range("E2").Value = "  +0.22%  "

$ws.Range("D3").Value = "1.799.79"
$ws.Range("E3").Value = "  +0.34%  "

$ws.Range("D4").Value = "'1.01"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.43%  "

$ws.Range("D5").Value = "'224.11"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.53%  "

$ws.Range("D6").Value = "'0.602"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.21%  "

$ws.Range("E7").Value = "  +0.37%  "

$ws.Range("D8").Value = "'41.18"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +13.90%  "

$ws.Range("E9").Value = "  -0.49%  "

$ws.Range("E10").Value = "  -1.52%  "

$ws.Range("D11").Value = "'0.0996"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.65%  "

$ws.Range("D12").Value = "2.058.38"
$ws.Range("E12").Value = "  +0.21%  "

$ws.Range("D13").Value = "1.795.14"
$ws.Range("E13").Value = "  -0.34%  "

$ws.Range("D14").Value = "'10.81"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.08%  "

$ws.Range("D15").Value = "34.387.19"
$ws.Range("E15").Value = "  +0.23%  "

$ws.Range("D16").Value = "'0.622"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.97%  "

$ws.Range("D17").Value = "'4.36"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.11%  "

$ws.Range("D18").Value = "'67.18"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.40%  "

$ws.Range("D19").Value = "'239.37"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.40%  "

$ws.Range("D20").Value = "0.0₃0763"
$ws.Range("E20").Value = "  -0.78%  "

$ws.Range("D21").Value = "'11.05"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.52%  "

$ws.Range("E22").Value = "  +0.32%  "

$ws.Range("D23").Value = "'4.08"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("E24").Value = "  -1.90%  "

$ws.Range("D25").Value = "'171.20"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.59%  "

$ws.Range("D26").Value = "'7.60"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.68%  "

$ws.Range("D27").Value = "'17.28"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.48%  "

$ws.Range("E28").Value = "  +0.16%  "

$ws.Range("D29").Value = "'1.01"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.58%  "

$ws.Range("D30").Value = "'1.22"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.46%  "

$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'3.82"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.24%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.0509"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.24%  "

$ws.Range("D34").Value = "'1.75"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.06%  "

$ws.Range("D35").Value = "1.314.88"
$ws.Range("E35").Value = "  -3.05%  "

$ws.Range("D36").Value = "'0.640"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.66%  "

$ws.Range("E37").Value = "  +0.56%  "

$ws.Range("D38").Value = "'84.75"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +5.34%  "

$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D39").Value = "'14.78"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +12.34%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.0187"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.08%  "

$ws.Range("E41").Value = "  +6.23%  "

$ws.Range("E42").Value = "  -0.68%  "

$ws.Range("E43").Value = "  +1.00%  "

$ws.Range("D44").Value = "'2.78"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("D45").Value = "'0.933"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.39%  "

$ws.Range("D46").Value = "'0.0518"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.63%  "

$ws.Range("D47").Value = "1.959.21"
$ws.Range("E47").Value = "  +0.21%  "

$ws.Range("D48").Value = "'5.81"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.33%  "

$ws.Range("E49").Value = "  +0.42%  "

$ws.Range("D50").Value = "'100.50"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.08%  "

$ws.Range("D51").Value = "'0.0610"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.41%  "
